$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 100
$ws_ALC.Range("H100").Value = 8001204.5
$ws_ALC.Range("I100").Value = 933.2353000000001
$ws_ALC.Range("J100").Value = 25001782
$ws_ALC.Range("K100").Value = 933.2353000000001
$ws_ALC.Range("L100").Value = 25001782
$ws_ALC.Range("M100").Value = -392.2353000000001
$ws_ALC.Range("N100").Value = -25002864

# Row 112
$ws_ALC.Range("H112").Value = 1226618.1
$ws_ALC.Range("J112").Value = 2316412
$ws_ALC.Range("L112").Value = 6949236
$ws_ALC.Range("N112").Value = -6951452

# Row 137
$ws_ALC.Range("H137").Value = 16282656
$ws_ALC.Range("I137").Value = 10417752
$ws_ALC.Range("K137").Value = 31253256
$ws_ALC.Range("M137").Value = -31250706


$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 32
$ws_ARM.Range("H32").Value = 6320101
$ws_ARM.Range("I32").Value = 1460307.5
$ws_ARM.Range("K32").Value = 1460307.5
$ws_ARM.Range("M32").Value = -1460020.5

# Row 61
$ws_ARM.Range("H61").Value = 1591018.1
$ws_ARM.Range("I61").Value = 898.61536
$ws_ARM.Range("J61").Value = 5349482.5
$ws_ARM.Range("K61").Value = 898.61536
$ws_ARM.Range("L61").Value = 5349482.5
$ws_ARM.Range("M61").Value = -686.61536
$ws_ARM.Range("N61").Value = -5349906.5

# Row 128
$ws_ARM.Range("H128").Value = 52490
$ws_ARM.Range("J128").Value = 52490
$ws_ARM.Range("L128").Value = 52490
$ws_ARM.Range("N128").Value = -62450

# Row 132
$ws_ARM.Range("H132").Value = 12157210
$ws_ARM.Range("I132").Value = 13894064
$ws_ARM.Range("J132").Value = 6946645.5
$ws_ARM.Range("K132").Value = 41682192
$ws_ARM.Range("L132").Value = 20839936.5
$ws_ARM.Range("M132").Value = -41679662
$ws_ARM.Range("N132").Value = -20844996.5

# Row 136
$ws_ARM.Range("H136").Value = 1591018.1
$ws_ARM.Range("I136").Value = 898.61536
$ws_ARM.Range("J136").Value = 5349482.5
$ws_ARM.Range("K136").Value = 2695.84608
$ws_ARM.Range("L136").Value = 16048447.5
$ws_ARM.Range("M136").Value = -145.8460800000003
$ws_ARM.Range("N136").Value = -16053547.5


$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 99
$ws_BSM.Range("H99").Value = 1218.5416
$ws_BSM.Range("I99").Value = 959.1875
$ws_BSM.Range("J99").Value = 1737.25
$ws_BSM.Range("K99").Value = 959.1875
$ws_BSM.Range("L99").Value = 1737.25
$ws_BSM.Range("M99").Value = 538.8125
$ws_BSM.Range("N99").Value = -4733.25

# Row 115
$ws_BSM.Range("H115").Value = 0
$ws_BSM.Range("J115").Value = 0
$ws_BSM.Range("L115").Value = 0
$ws_BSM.Range("N115").ClearContents()

# Row 134
$ws_BSM.Range("H134").Value = 12459561
$ws_BSM.Range("I134").Value = 17242278
$ws_BSM.Range("J134").Value = 2552503.5
$ws_BSM.Range("K134").Value = 51726834
$ws_BSM.Range("L134").Value = 7657510.5
$ws_BSM.Range("M134").Value = -51724299
$ws_BSM.Range("N134").Value = -7662580.5


$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 31
$ws_CRP.Range("H31").Value = 2234917.2
$ws_CRP.Range("I31").Value = 2472.923
$ws_CRP.Range("J31").Value = 4169702.2
$ws_CRP.Range("K31").Value = 2472.923
$ws_CRP.Range("L31").Value = 4169702.2
$ws_CRP.Range("M31").Value = -2177.923
$ws_CRP.Range("N31").Value = -4170292.2

# Row 34
$ws_CRP.Range("H34").Value = 2234917.2
$ws_CRP.Range("I34").Value = 2472.923
$ws_CRP.Range("J34").Value = 4169702.2
$ws_CRP.Range("K34").Value = 2472.923
$ws_CRP.Range("L34").Value = 4169702.2
$ws_CRP.Range("M34").Value = -2270.923
$ws_CRP.Range("N34").Value = -4170106.2

# Row 99
$ws_CRP.Range("H99").Value = 20699.908
$ws_CRP.Range("I99").Value = 16259.8
$ws_CRP.Range("K99").Value = 16259.8
$ws_CRP.Range("M99").Value = -14761.8

# Row 107
$ws_CRP.Range("H107").Value = 802.5263
$ws_CRP.Range("I107").Value = 278.5
$ws_CRP.Range("J107").Value = 1183.6364
$ws_CRP.Range("K107").Value = 278.5
$ws_CRP.Range("L107").Value = 1183.6364
$ws_CRP.Range("M107").Value = 1641.5
$ws_CRP.Range("N107").Value = -5023.6364

# Row 126
$ws_CRP.Range("H126").Value = 20699.908
$ws_CRP.Range("I126").Value = 16259.8
$ws_CRP.Range("K126").Value = 48779.39999999999
$ws_CRP.Range("M126").Value = -46309.39999999999

# Row 131
$ws_CRP.Range("H131").Value = 15793.2
$ws_CRP.Range("J131").Value = 15793.2
$ws_CRP.Range("L131").Value = 15793.2
$ws_CRP.Range("N131").Value = -25873.2

# Row 132
$ws_CRP.Range("H132").Value = 2429.2334
$ws_CRP.Range("I132").Value = 1777.125
$ws_CRP.Range("J132").Value = 5037.6665
$ws_CRP.Range("K132").Value = 5331.375
$ws_CRP.Range("L132").Value = 15112.9995
$ws_CRP.Range("M132").Value = -2801.375
$ws_CRP.Range("N132").Value = -20172.9995

# Row 134
$ws_CRP.Range("H134").Value = 5525.2
$ws_CRP.Range("I134").Value = 5229.68
$ws_CRP.Range("J134").Value = 7002.8
$ws_CRP.Range("K134").Value = 15689.04
$ws_CRP.Range("L134").Value = 21008.4
$ws_CRP.Range("M134").Value = -13154.04
$ws_CRP.Range("N134").Value = -26078.4


$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 131
$ws_CUL.Range("H131").Value = 22537.238
$ws_CUL.Range("I131").Value = 111488.78
$ws_CUL.Range("J131").Value = 900.37836
$ws_CUL.Range("K131").Value = 334466.34
$ws_CUL.Range("L131").Value = 2701.13508
$ws_CUL.Range("M131").Value = -329426.34
$ws_CUL.Range("N131").Value = -12781.13508


$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 27
$ws_GSM.Range("H27").Value = 0
$ws_GSM.Range("J27").Value = 0
$ws_GSM.Range("L27").Value = 0
$ws_GSM.Range("N27").ClearContents()

# Row 132
$ws_GSM.Range("H132").Value = 16668430
$ws_GSM.Range("I132").Value = 40000756
$ws_GSM.Range("J132").Value = 2484.5715
$ws_GSM.Range("K132").Value = 120002268
$ws_GSM.Range("L132").Value = 7453.7145
$ws_GSM.Range("M132").Value = -119999738
$ws_GSM.Range("N132").Value = -12513.7145

# Row 139
$ws_GSM.Range("H139").Value = 54125
$ws_GSM.Range("J139").Value = 54125
$ws_GSM.Range("L139").Value = 54125
$ws_GSM.Range("N139").Value = -64405


$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 7
$ws_LTW.Range("H7").Value = 1296.2
$ws_LTW.Range("I7").Value = 1273.5555
$ws_LTW.Range("J7").Value = 1500
$ws_LTW.Range("K7").Value = 1273.5555
$ws_LTW.Range("L7").Value = 1500
$ws_LTW.Range("M7").Value = -1161.5555
$ws_LTW.Range("N7").Value = -1724

# Row 30
$ws_LTW.Range("H30").Value = 0
$ws_LTW.Range("I30").Value = 0
$ws_LTW.Range("K30").Value = 0
$ws_LTW.Range("M30").ClearContents()

# Row 126
$ws_LTW.Range("H126").Value = 1296.2
$ws_LTW.Range("I126").Value = 1273.5555
$ws_LTW.Range("J126").Value = 1500
$ws_LTW.Range("K126").Value = 3820.6665
$ws_LTW.Range("L126").Value = 4500
$ws_LTW.Range("M126").Value = -1350.6665
$ws_LTW.Range("N126").Value = -9440

# Row 132
$ws_LTW.Range("H132").Value = 5499789
$ws_LTW.Range("I132").Value = 7148818.5
$ws_LTW.Range("K132").Value = 21446455.5
$ws_LTW.Range("M132").Value = -21443925.5


$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 107
$ws_WVR.Range("H107").Value = 15407.9
$ws_WVR.Range("I107").Value = 25035
$ws_WVR.Range("J107").Value = 8989.833000000001
$ws_WVR.Range("K107").Value = 75105
$ws_WVR.Range("L107").Value = 26969.499
$ws_WVR.Range("M107").Value = -73185
$ws_WVR.Range("N107").Value = -30809.499

# Row 132
$ws_WVR.Range("H132").Value = 2028.4222
$ws_WVR.Range("I132").Value = 2006.2778
$ws_WVR.Range("J132").Value = 2117
$ws_WVR.Range("K132").Value = 6018.8334
$ws_WVR.Range("L132").Value = 6351
$ws_WVR.Range("M132").Value = -3488.8334
$ws_WVR.Range("N132").Value = -11411

# Row 138
$ws_WVR.Range("H138").Value = 48456.332
$ws_WVR.Range("J138").Value = 48456.332
$ws_WVR.Range("L138").Value = 48456.332
$ws_WVR.Range("N138").Value = -58736.332

